$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the input values that drive the MMULT array formula (E24:E27)
# and the equivalent SUMPRODUCT-style formulas (H24:H27).
$ws.Range("B24").Value = 0.1
$ws.Range("B26").Value = -0.1

# Move the active selection to H24:H27 (active cell H24), matching the
# author's saved view state.
$ws.Range("H24:H27").Select()

$wb.Save()
